# "Generate Report for Archive"
# - Status moves from "Ready for handoff" to "In Translation" everywhere it
#   appears (Overview!E2/F2 and the per-language Status column, zh-cn!C2 /
#   de-de!C2, all share the one "Ready for handoff" shared string).
# - The Status columns (Overview E:F, zh-cn C, de-de C) get narrower to fit
#   the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns(5).ColumnWidth = 12.5
$overview.Columns(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns(3).ColumnWidth = 12.5
